# Water Data SA Download Code
# Adds an "isSonde" flag column and a "Conversion" column to the Vars sheet,
# renames the LEVEL variable's AED name to "H", zeroes out the old
# Conversion values (now in column C / isSonde), and appends a new
# fDOM_SONDE / fDOM row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("C1").Value = "isSonde"
$ws.Range("D1").Value = "Conversion"

# --- Rename the LEVEL row's AED name to "H" ----------------------------
$ws.Range("B6").Value = "H"

# --- New "isSonde" flag in column C: 0 for non-sonde rows, keep 1 for
#     rows that already represent *_SONDE variables ---------------------
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0
# C11 (TEMP_SONDE) stays 1
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0
# C16 (WQ_CAR_PH_SONDE) stays 1
# C17 (WQ_DIAG_TOT_TCHLA_SONDE) stays 1
$ws.Range("C18").Value = 0
# C19 (WQ_DIAG_TOT_TURBIDITY_SONDE) stays 1
# C20 (WQ_OXY_OXY_SAT_SONDE) stays 1
# C21 (WQ_OXY_OXY_SONDE) stays 1
# C22 (WQ_PHY_BGA_SONDE) stays 1

# --- New "Conversion" values/formulas in column D ----------------------
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Formula = "=1000/86400"
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("D21").Formula = "=1000/32"
$ws.Range("D22").Value = 1

# --- New fDOM_SONDE / fDOM row ------------------------------------------
$ws.Range("A23").Value = "fDOM_SONDE"
$ws.Range("B23").Value = "fDOM"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1

# --- Set the new column widths / active selection to match the edit ----
$ws.Columns.Item(3).ColumnWidth = 21.7265625
$ws.Range("E23").Select()
